# TC_55: rename sheet, update labels, expand quarterly series, restyle number formats,
# and refresh the CEIC add-in comment payload on A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet tab.
$ws.Name = "Data"

# 2) Update the header / label text that changed.
$ws.Range("C1").Value = "Real Estate Investment: Residential: Hebei [DISAGGREGATE(Quarterly; March, June, September, December; Replicate)]"
$ws.Range("A11").Value = "Function Information"

# 3) Rewrite row 13 (first data row) in place, then append the rest of the
#    disaggregated quarterly series in rows 14-48. Number formats are first
#    applied using the *original* format strings so new cells land on the
#    same style slots (numFmtId 164 / 165) as the pre-existing rows, then
#    every cell sharing those formats is re-formatted to the new codes in
#    one sweep below.
$rows = @(
    ,@(13, 35490, $null, 3271.34)
    ,@(14, 35582, $null, 3271.34)
    ,@(15, 35674, $null, 3271.34)
    ,@(16, 35765, 3271.34, 3271.34)
    ,@(17, 36220, $null, 6384.48)
    ,@(18, 36312, $null, 6384.48)
    ,@(19, 36404, $null, 6384.48)
    ,@(20, 36495, 6384.48, 6384.48)
    ,@(21, 36586, $null, 7111.74)
    ,@(22, 36678, $null, 7111.74)
    ,@(23, 36770, $null, 7111.74)
    ,@(24, 36861, 7111.74, 7111.74)
    ,@(25, 36951, $null, 8354.02)
    ,@(26, 37043, $null, 8354.02)
    ,@(27, 37135, $null, 8354.02)
    ,@(28, 37226, 8354.02, 8354.02)
    ,@(29, 37316, $null, 10540.67)
    ,@(30, 37408, $null, 10540.67)
    ,@(31, 37500, $null, 10540.67)
    ,@(32, 37591, 10540.67, 10540.67)
    ,@(33, 37681, $null, 16421.36)
    ,@(34, 37773, $null, 16421.36)
    ,@(35, 37865, $null, 16421.36)
    ,@(36, 37956, 16421.36, 16421.36)
    ,@(37, 38047, $null, 22352.63)
    ,@(38, 38139, $null, 22352.63)
    ,@(39, 38231, $null, 22352.63)
    ,@(40, 38322, 22352.63, 22352.63)
    ,@(41, 38412, $null, 29205.22)
    ,@(42, 38504, $null, 29205.22)
    ,@(43, 38596, $null, 29205.22)
    ,@(44, 38687, 29205.22, 29205.22)
    ,@(45, 38777, $null, 37962.97)
    ,@(46, 38869, $null, 37962.97)
    ,@(47, 38961, $null, 37962.97)
    ,@(48, 39052, 37962.97, 37962.97)
)

foreach ($row in $rows) {
    $r = $row[0]
    $dateVal = $row[1]
    $bVal = $row[2]
    $cVal = $row[3]

    $aCell = $ws.Cells.Item($r, 1)
    $aCell.NumberFormat = "MM/yyyy"
    $aCell.Value = $dateVal

    $bCell = $ws.Cells.Item($r, 2)
    $bCell.NumberFormat = "0.000"
    if ($bVal -ne $null) {
        $bCell.Value = $bVal
    }

    $cCell = $ws.Cells.Item($r, 3)
    $cCell.NumberFormat = "0.000"
    $cCell.Value = $cVal
}

# 4) Re-point every cell that used the old "MM/yyyy" / "0.000" custom
#    formats (the First Obs. Date row plus the whole data block) at the
#    new format codes.
$ws.Range("B12:C12").NumberFormat = "dd/MM/yyyy"
$ws.Range("A13:A48").NumberFormat = "dd/MM/yyyy"
$ws.Range("B13:C48").NumberFormat = "###0.000"

# 5) Refresh the CEIC add-in metadata payload stashed in the A1 cell comment.
$comment = $ws.Range("A1").Comment
$newBlob = "yjMAAB+LCAAAAAAAAAPtWltvG8cV/isLPiWAqdklZd0y3oAiKYWtKMokFVspimC5OxS3Xu4ye5HEtxRokSJNURSFU6RX9ClFgbpGmwCp3ct/CSzZeepf6JmZvcxeKHFlBW0DG4a9c24zc2bmnO8MB795NrWkE+J6pmPfqSgrckUitu4Ypn18pxL446qyVnlTxe0znVgHmqtNiQ/CEmjZ3taZZ96pTHx/toXQ6enpyml9xXGPUU2WFXS/uzfQJ2SqVU3b8zVbJ5VYy7haq6LipjHtEl8zNF/jmncqnUFnpUlMvQW0rmZrx8Rd2Q480yae17Z90zeJRzVdovmk2eq+zSem1lbWVhSMcvREcjswLYPLpSQ5PZSDbsnQnBK1JisbVXmjWleGirylKFv1tZXV2sY7kWIsiPc0zx8Q98TUGWHga9MZU5c36ooiK0p9DaNCIbCVOEDFPcvokxPTI0aTWJZXyiMoXMCG7sOsyzlTxkjQDQ1dfwi7rjabDE3fIuWG0e9uS1M7HEtiRMU7jkt08N+1hrRPTntu6NbhbA+4w4np+vOWNi9t69Ajbm9GnVROVcUtx/YbFnH9wxmsNTFgKwBD9d2AYLSAmSi1TE+Hb9MOiKGONcsTlVJMfM9xH3gzTSf7cI4RtXFqW45mwIbzTc83dS8xkOPgA9eZgUnofduxjB0wGw6xgBGb7tjgY9rvtuM8yFpPMzFbVrbAsKhTzY/Ec3Q8mDinPduaD4KRp7vmiBit7Ui6kIfpiQy1m4HnO1MYRULCnCZQDAN1u2gOf+AcZpm4RXRzqlkHFvjSU+tgK0XAjcB3xqbfdKxgasdOzVDxPZjXkJzF84zbuAdrbFPXO3bHjuS5swtZaYW+cxr3mWcwVwjkhqdHWy3PyAq3gBYtYp7D1oXOcse0IE2IKyJQ03tjMCHEL9wYnINpRNyhiUfdntM+MUooGPYnbHKgqsrm5npVVuDvUJa32F/oOWbjtm2wD0gza1WlRkN3Ihcx8X4w7Y3gIJ+wOakK8DIkDLOwti3NfgDUe6Y/2W9Eoy/gYD7nhfJ5HoYjO7O0OSPHfhFpuGPrVmAQHgk69phtSjo2vowL2ThH2oOzrWLNng/nMwjInrnlw8edCqToLc93AQRUVN0JbN+d05CBUSh6lY4XjGzWgWYtrTN2yXsBYI/5TmDrTcdYvjeDe+fQNv3lR+gELo+Dy6sw79GIGHgtQmMLi/ZL6+tl5uS5pcSnNpk6tqkv721wMh29cY2JeNGpEjTQgr3VImMtsADd+JCEjpP8kiHjhvcgKyOS8KFrRcFBpdjRA/CoG9MVHdIrBUgrujOlBASY7d4AI1GeYgSdtO3jPc0+DiALxwcwS49DE00gQ1ezPTqdOOlmolSxEI4ONAcDKj/lvYB5jJ9yB7gYZeTwkExnjqtZXXCMuROuT4goIFl3NX8StiDsW0SPnIwS1VgrPbJo4FeJsfjNp0FPRhhPMkQmROfCYWoik9AwnWUX9q/V1Cxz5PLwE2W5Ih4sWAKfokBFJ1cSSkVrAKUKJKZvkznFr0kjpLMtq0QMvoFpxFEH/dWN2m25XoOkT9uYzbhPNEtqw673idSxT4jnT0FtS+oTzzTgy9SsLektMiImZAvmojBel9YW9fBOFBDZUBow3jQlLQCp+NiEeJsXjDmJgnpENNeaC4J8qnuODnIXP/rX+a8fP3vy6cWHD198/oN///2Xz/7x0/NHP4SPi7/89fyjX/BpcmE81EYWYQMabm9syPVV2GcxCVPnIoYdjUD3Ge3oiEHGuI3D0oc1mu1Oc3dvm8WTmBip89iLaFU1d4KkOeCTYB2xJUXRTuAi6jCKT2E7xRViuUqrnBOSlhb5ixS5L54//fT50z8t1A4dlsIstykWuRqzKEpOLsYstHrkxUFUm65W5dvVWk0QzsjgPpTYUAvEfuoYal2RN+VaXVbiWG7EG7lIKMsKLQ21Y5TR46QmxxHxFhDbEZNt/CEckZjNj4LQCLfo5z9+8eeHKanQuyElbQUGx/I97QxFDWZ6vz+UBr3DfrMtDdsDuk8SniDHjV8iHPYen6fUprLtQLNuScSm1ypSBYqFiuSMJaLpE2kOJ1E4h6nNVkTlHV3TZHaUu64TzPiKCAoJtUAyjiaFGgWxhvGYP3NBJ2EViPOxnv/tsyKFcCKtBPnFNwQiDac4nCTww1P7yT+fffHBsydPLh7/7PyL76cshP3EhTPsczhNYjPe9hDywnyToeB7A+bMB/K7Qn4JibTiOHBM2/dUZY0VG2ELg6pCrbH/cWcKKY8ZZv4CeoaC39K89pkfHmx1H6M0AcY50yDbOklZFhN4DE/8+tVvfnvxq88uPn784oM/nn/4h/OPPn7+9HcvHv2en7qLh48vfvIojPLZRMDGQos9DgIldoGgS/Q0SjR3S1++/3PJdnwJIIcUsIj05fufCMboQBk4SSwDpIsHkh5CTlRUpnqSMJR4DCm9WIUDgCZNYfVYIkxizszUk07eqVJT9NwxxmudYTXwiOQAmnodZpIWTpSX1QtVeEo9WJdrSi3k8tHQKYw0T3D9ruWMAGREDFabZ0RSWpcrJLKsv9293nZjLxHhg+i5BpTyMr1goB84gpQ0pXS8qBVtNYECXAB+emDR65ScWJ4VWxbCGArvJcYNg4Y/VZZZ4Z/OjykJ3AxclwMiO7ztHgQzAMPRHdZiPrvYE/DvPseqIiJO2p1Wmg9tgQuJMM2mBMZnoSlk8TDV8egVCIez+9Q1SRN4qctAcEd4oc2R1gngShfRuNN2XcctDD4JJxLrApKGiIISj8cybE056jaStYoIUcD7r6Nv6TutzqCxu9tv7zaG7dfuBgCkCGScN6Su5uqTW9K3ApvckiDu+WQ6Iu4tCdIm+3oDrM0sU4c+Xv/uy4H4GxvEJbXA3WVrgZxgUS0QD/G65cCN+v1VVXEjVYVcf1VVfLOriviYJVUAHDfhsLFjVilfW7yU4cUVxt3SFUahxqUVRlEo+yYWGcriIkN511BuK4pe36hCLh1XV8ejtermOiFVeVUfrdc2tZq2ullUf6ytfq31R7RwBb9LLBzD8mnvmnvgsp1ZsKw3HWRv5J7w8lh7RQy8unpFlyzcq5ryVU150zWlkq8p+Y8fmZ9EICYZagQ8DyAeUWAhZLibwKT0hyohZyYPoPiI4nZYwIzd96LS5W3NCiCyHDb6w3Z/D4AJJzDkqjZ7+4NhYx9SffjjWcZYkW1IxQdwAB0j3QNdxzKmkTgHDmuztZdAjZcgCYsN1/QnYMDUezPCf73xIH7m5GJVWuMmSlKRFiuDszcKYRxyuAK5blWJeCmDUsEWZVJLth2FPJGOhJoICQUKKipHUKr+QC9da6CFNcJiTjTSJCsV3pMUJK1CuUWVQaHw1ZVBFs4zyvJgH6UyG8qD+YSEFuJ0lAPjqBgYLyDnVOLtX0gs2lUxMikkRh0IYBFxYILyGBPdzD01EvGY/FKQcHl4gAQCyqfrNAml8ivK5lKRgITEhwrSXIaG0nkKhVEJ/T/fb4Z9/y9fcKLrXGmim7u5FJJMyxyPjwBH3k+eXkUEPuLECm3w/JpDJkiALsLeCN9n3exmuSGgU2T51S172Vv25avvr+UCHiUPZxIA1SIW4L2yz5kj7a5zcm1d2BZlVTtezzJCP5d7+xO7JTEgvjmne+hln5zzfdhwXYAl9HFq6Tfi0cuxvmYflxwNnwpTpE+woPcINZiu59+n6Tf84pSjmHLEIdt9tcZx2X3ePlJX1zkBBJBoHaWGGZ1qn7/8d6w9c2qWfJclR0c/bQR8OZtxBNYpt1NoHb5PznyMBAuQrUffA2DAny6WscY3LJQesT59ZeuZxxO/7MDWRxoxyEiu6iNSq64a8kZ1k5B6VVHgX02v1WT5Nn2jGxqHoGKS05KdoGjBhGLuPydh6tDKMwAA"
$comment.Text($newBlob)
